$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.225.63"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "1.645.21"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "

$ws.Range("E11").Value = "  +0.02%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "1.873.32"
$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").Value = "1.643.84"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.45%  "

$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "26.212.26"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "

$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.50"
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.35%  "

$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0504"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  +1.42%  "

$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.913"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("E37").Value = "  +1.74%  "

$ws.Range("D38").Value = "1.134.25"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("E39").Value = "  -1.66%  "

$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("D45").Value = "1.782.59"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.56%  "

$ws.Range("E47").Value = "  +4.10%  "

$ws.Range("E48").Value = "  +2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.10%  "

$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("E51").Value = "  +1.52%  "
